# Add new columns I (I0) and J (IF) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Headers - match formatting of existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for rows 2-70: column I (I0), column J (IF)
$data = @(
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(6,7),
    @(4,4),
    @(7,7),
    @(8,8),
    @(6,7),
    @(7,8),
    @(7,8),
    @(9,9),
    @(6,6),
    @(8,8),
    @(6,6),
    @(6,7),
    @(7,7),
    @(10,10),
    @(8,9),
    @(8,9),
    @(4,4),
    @(8,8),
    @(7,7),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(5,5),
    @(6,7),
    @(5,5),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,7),
    @(9,9),
    @(6,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(7,7),
    @(7,7),
    @(6,6),
    @(7,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,8),
    @(10,10),
    @(6,7),
    @(10,10),
    @(6,6),
    @(7,8),
    @(7,7),
    @(11,11),
    @(7,8),
    @(7,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(5,5),
    @(3,3),
    @(5,5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}

$wb.Save()
